# The author added a new weekly price record for "Feria Lagunitas de Puerto
# Montt - Ciboulette" (commit: "Fruta / hortaliza, semanal"). In the
# consolidated log, newest entries are inserted right after the header block
# of rows (at row 51), pushing the previously-existing rows 51-128 down by
# one (to 52-129), and the sheet's used range grows from A1:R128 to A1:R129.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 51; everything at/after row 51 shifts down.
$ws.Rows.Item(51).Insert()

# Populate the new row 51 with the new daily record.
$ws.Range("A51").Value = 4
$ws.Range("B51").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C51").Value = "Los Lagos"
$ws.Range("D51").Value = 44495
$ws.Range("E51").Value = 10
$ws.Range("F51").Value = 100112039
$ws.Range("G51").Value = "Ciboulette"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 240
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3000
$ws.Range("M51").Value = 2750
$ws.Range("N51").Value = "$/docena de atados"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 917
$ws.Range("Q51").Value = 3
$ws.Range("R51").Value = "Hortaliza"
